# Fix the DaPaiWei deceased date issue pointed out by Brother Guo and Sister Chan.
#  - Format check
#  - Valid date range check
#
# Sheet 6 ("12個月內往生親友牌位") gets 3 extra sample rows so the format /
# date-range checks can be exercised, and becomes the active sheet/tab.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(6)

# A new row is inserted right under the existing single data row (old row 4),
# pushing that row down to row 5 and inheriting row 3's formatting for the
# freshly inserted row 4.
$ws.Rows.Item(4).Insert()

# Two more rows are appended at the bottom of the table (rows 6 and 7).
$ws.Cells.Item(7, 2).Value = "CCCC"
$ws.Cells.Item(7, 3).Value = "2018-09-30"
$ws.Cells.Item(7, 5).Value = "DDDD"

$ws.Cells.Item(6, 2).Value = "MMMM"
$ws.Cells.Item(6, 3).Value = "2019-02-28"
$ws.Cells.Item(6, 5).Value = "DDDD"

$ws.Cells.Item(4, 2).Value = "ZZZZ"
$ws.Cells.Item(4, 3).Value = "2018-07-01"
$ws.Cells.Item(4, 4).Value = "UUUU"
$ws.Cells.Item(4, 5).Value = "VVVV"

$ws.Cells.Item(4, 1).Value = "XXXX"

$ws.Cells.Item(6, 1).Value = "朋有"
$ws.Cells.Item(6, 4).Value = "朋有"

$ws.Cells.Item(7, 1).Value = "朋有"
$ws.Cells.Item(7, 4).Value = "朋有"

# Give the two appended rows (6 & 7) the same direct formatting as row 5
# (the row they were appended below), matching the formatting row 4 gets
# automatically from the Insert() above.
$ws.Range("A5:E5").Copy()
$ws.Range("A6:E7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# The sheet-scoped named range used for data validation on this table needs
# to grow by the 3 rows that were just added.
$name = $ws.Names.Item(1)
$name.RefersTo = "='12個月內往生親友牌位'!`$B`$3:`$F`$522"

$ws.Activate()
$ws.Range("F6").Select()
